$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a numeric-looking string that must stay text
# (so formatting such as trailing zeros, e.g. "122.80", is preserved exactly
# like the other already-textual price cells in this sheet).
$textCells = @('D4', 'D5', 'D7', 'D8', 'D9', 'D11', 'D12', 'D13', 'D14', 'D15', 'D17', 'D18', 'D19', 'D21', 'D22', 'D24', 'D25', 'D26', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D41', 'D42', 'D43', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '28.600.46'
$ws.Range('E2').Value = '  +0.78%  '
$ws.Range('D3').Value = '1.802.42'
$ws.Range('E3').Value = '  -0.62%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '316.85'
$ws.Range('E5').Value = '  -0.32%  '
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('D7').Value = '0.5419'
$ws.Range('E7').Value = '  -5.37%  '
$ws.Range('D8').Value = '0.3766'
$ws.Range('E8').Value = '  -2.84%  '
$ws.Range('D9').Value = '0.07488'
$ws.Range('E9').Value = '  -1.67%  '
$ws.Range('E10').Value = '  -1.80%  '
$ws.Range('D11').Value = '1.114'
$ws.Range('E11').Value = '  -2.27%  '
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  +0.11%  '
$ws.Range('D13').Value = '20.64'
$ws.Range('D14').Value = '6.147'
$ws.Range('E14').Value = '  -1.79%  '
$ws.Range('D15').Value = '7.387'
$ws.Range('E15').Value = '  +1.49%  '
$ws.Range('D16').Value = '1.797.52'
$ws.Range('E16').Value = '  -0.92%  '
$ws.Range('D17').Value = '90.21'
$ws.Range('D18').Value = '0.00001065'
$ws.Range('E18').Value = '  -1.12%  '
$ws.Range('D19').Value = '0.06446'
$ws.Range('E19').Value = '  -0.49%  '
$ws.Range('E20').Value = '  +0.08%  '
$ws.Range('D21').Value = '17.23'
$ws.Range('D22').Value = '5.921'
$ws.Range('E22').Value = '  -1.32%  '
$ws.Range('D23').Value = '28.628.05'
$ws.Range('E23').Value = '  +0.81%  '
$ws.Range('D24').Value = '11.11'
$ws.Range('E24').Value = '  -1.77%  '
$ws.Range('D25').Value = '2.092'
$ws.Range('E25').Value = '  -1.59%  '
$ws.Range('D26').Value = '158.55'
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('E27').Value = '  -2.29%  '
$ws.Range('D28').Value = '2.008.10'
$ws.Range('E28').Value = '  -0.73%  '
$ws.Range('D29').Value = '2.349'
$ws.Range('E29').Value = '  -3.84%  '
$ws.Range('D30').Value = '122.80'
$ws.Range('E30').Value = '  -1.18%  '
$ws.Range('D31').Value = '1.104'
$ws.Range('E31').Value = '  -5.43%  '
$ws.Range('D32').Value = '0.1057'
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('D33').Value = '5.638'
$ws.Range('E33').Value = '  -2.62%  '
$ws.Range('D34').Value = '3.681'
$ws.Range('E34').Value = '  +1.30%  '
$ws.Range('D35').Value = '0.06492'
$ws.Range('E35').Value = '  +6.50%  '
$ws.Range('D36').Value = '0.2247'
$ws.Range('E36').Value = '  +3.74%  '
$ws.Range('D37').Value = '0.02299'
$ws.Range('E37').Value = '  -0.89%  '
$ws.Range('D38').Value = '8.734'
$ws.Range('E38').Value = '  -1.96%  '
$ws.Range('D39').Value = '5.021'
$ws.Range('E39').Value = '  -0.38%  '
$ws.Range('E40').Value = '  -3.83%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = '1.204'
$ws.Range('E41').Value = '  +3.51%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = '0.6221'
$ws.Range('E42').Value = '  -3.27%  '
$ws.Range('D43').Value = '1.433'
$ws.Range('E43').Value = '  +4.22%  '
$ws.Range('E44').Value = '  +0.10%  '
$ws.Range('D45').Value = '13.22'
$ws.Range('E45').Value = '  -2.03%  '
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').Value = '3.689'
$ws.Range('E46').Value = '  -0.59%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '0.5842'
$ws.Range('E47').Value = '  -2.82%  '
$ws.Range('D48').Value = '126.48'
$ws.Range('E48').Value = '  +2.95%  '
$ws.Range('D49').Value = '1.938'
$ws.Range('E49').Value = '  -0.40%  '
$ws.Range('D50').Value = '1.155'
$ws.Range('E50').Value = '  +0.36%  '
$ws.Range('D51').Value = '0.06887'
$ws.Range('E51').Value = '  +0.45%  '
